$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily data rows 270-301 (dates 2021-05-28 through 2021-06-28),
# columns: row, A (date serial), B (nuovi pos.), C (somma mobile 7gg.), D (somma mobile 7gg. per 100mila abitanti)
$newRows = @(
    @(270, 44344, 0, 0, 0),
    @(271, 44345, 0, 0, 0),
    @(272, 44346, 0, 0, 0),
    @(273, 44347, 0, 0, 0),
    @(274, 44348, 0, 0, 0),
    @(275, 44349, 0, 0, 0),
    @(276, 44350, 0, 0, 0),
    @(277, 44351, 0, 0, 0),
    @(278, 44352, 0, 0, 0),
    @(279, 44353, 0, 0, 0),
    @(280, 44354, 0, 0, 0),
    @(281, 44355, 0, 0, 0),
    @(282, 44356, 0, 0, 0),
    @(283, 44357, 0, 0, 0),
    @(284, 44358, 1, 1, 25.4323499491353),
    @(285, 44359, 0, 1, 25.4323499491353),
    @(286, 44360, 0, 1, 25.4323499491353),
    @(287, 44361, 0, 1, 25.4323499491353),
    @(288, 44362, 0, 1, 25.4323499491353),
    @(289, 44363, 0, 1, 25.4323499491353),
    @(290, 44364, 0, 1, 25.4323499491353),
    @(291, 44365, 0, 0, 0),
    @(292, 44366, 1, 1, 25.4323499491353),
    @(293, 44367, 0, 1, 25.4323499491353),
    @(294, 44368, 1, 2, 50.8646998982706),
    @(295, 44369, 0, 2, 50.8646998982706),
    @(296, 44370, 0, 2, 50.8646998982706),
    @(297, 44371, 0, 2, 50.8646998982706),
    @(298, 44372, 0, 2, 50.8646998982706),
    @(299, 44373, 0, 1, 25.4323499491353),
    @(300, 44374, 0, 1, 25.4323499491353),
    @(301, 44375, 0, 0, 0)
)

foreach ($r in $newRows) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 1).Value = $r[1]
    $ws.Cells.Item($rowNum, 2).Value = $r[2]
    $ws.Cells.Item($rowNum, 3).Value = $r[3]
    $ws.Cells.Item($rowNum, 4).Value = $r[4]
}

# Copy the date-column formatting (style used on existing rows) onto the new A-column cells
$ws.Range("A269").Copy()
$ws.Range("A270:A301").PasteSpecial(-4122)
$excel.CutCopyMode = 0
